$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5997.8335
$ws.Range("J18").Value = 15999
$ws.Range("L18").Value = 15999
$ws.Range("N18").Value = -16567
$ws.Range("H62").Value = 13596.173
$ws.Range("I62").Value = 11723.177
$ws.Range("K62").Value = 11723.177
$ws.Range("M62").Value = -11099.177
$ws.Range("H65").Value = 13596.173
$ws.Range("I65").Value = 11723.177
$ws.Range("K65").Value = 58615.88499999999
$ws.Range("M65").Value = -55495.88499999999
$ws.Range("H69").Value = 7917.5557
$ws.Range("I69").Value = 5600
$ws.Range("J69").Value = 8207.25
$ws.Range("K69").Value = 16800
$ws.Range("L69").Value = 24621.75
$ws.Range("M69").Value = -15926
$ws.Range("N69").Value = -26369.75
$ws.Range("H72").Value = 7917.5557
$ws.Range("I72").Value = 5600
$ws.Range("J72").Value = 8207.25
$ws.Range("K72").Value = 50400
$ws.Range("L72").Value = 73865.25
$ws.Range("M72").Value = -46032
$ws.Range("N72").Value = -82601.25
$ws.Range("H96").Value = 1158.8235
$ws.Range("I96").Value = 726.75
$ws.Range("K96").Value = 2180.25
$ws.Range("M96").Value = -807.25
$ws.Range("H100").Value = 5253.2383
$ws.Range("I100").Value = 1917.7273
$ws.Range("J100").Value = 8922.299999999999
$ws.Range("K100").Value = 1917.7273
$ws.Range("L100").Value = 8922.299999999999
$ws.Range("M100").Value = -1376.7273
$ws.Range("N100").Value = -10004.3

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20419.768
$ws.Range("I32").Value = 21164.166
$ws.Range("K32").Value = 21164.166
$ws.Range("M32").Value = -20877.166
$ws.Range("H61").Value = 2779197.2
$ws.Range("I61").Value = 5556794.5
$ws.Range("J61").Value = 1599.8334
$ws.Range("K61").Value = 5556794.5
$ws.Range("L61").Value = 1599.8334
$ws.Range("M61").Value = -5556582.5
$ws.Range("N61").Value = -2023.8334
$ws.Range("H74").Value = 4307.5454
$ws.Range("I74").Value = 1230.5
$ws.Range("K74").Value = 1230.5
$ws.Range("M74").Value = -356.5
$ws.Range("H77").Value = 4307.5454
$ws.Range("I77").Value = 1230.5
$ws.Range("K77").Value = 6152.5
$ws.Range("M77").Value = -1784.5
$ws.Range("H122").Value = 2367.3667
$ws.Range("I122").Value = 2149
$ws.Range("K122").Value = 6447
$ws.Range("M122").Value = -3997
$ws.Range("H136").Value = 2779197.2
$ws.Range("I136").Value = 5556794.5
$ws.Range("J136").Value = 1599.8334
$ws.Range("K136").Value = 16670383.5
$ws.Range("L136").Value = 4799.5002
$ws.Range("M136").Value = -16667833.5
$ws.Range("N136").Value = -9899.5002
$ws.Range("H137").Value = 37889.75
$ws.Range("I137").Value = 25779
$ws.Range("K137").Value = 25779
$ws.Range("M137").Value = -20679

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 72000
$ws.Range("J55").Value = 72000
$ws.Range("L55").Value = 72000
$ws.Range("N55").Value = -72546
$ws.Range("H107").Value = 2708.3635
$ws.Range("I107").Value = 2121
$ws.Range("J107").Value = 3736.25
$ws.Range("K107").Value = 2121
$ws.Range("L107").Value = 3736.25
$ws.Range("M107").Value = -201
$ws.Range("N107").Value = -7576.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null
$ws.Range("H31").Value = 18140.074
$ws.Range("I31").Value = 5672.1904
$ws.Range("J31").Value = 61777.668
$ws.Range("K31").Value = 5672.1904
$ws.Range("L31").Value = 61777.668
$ws.Range("M31").Value = -5377.1904
$ws.Range("N31").Value = -62367.668
$ws.Range("H34").Value = 18140.074
$ws.Range("I34").Value = 5672.1904
$ws.Range("J34").Value = 61777.668
$ws.Range("K34").Value = 5672.1904
$ws.Range("L34").Value = 61777.668
$ws.Range("M34").Value = -5470.1904
$ws.Range("N34").Value = -62181.668
$ws.Range("H50").Value = 46665.332
$ws.Range("J50").Value = 59998
$ws.Range("L50").Value = 59998
$ws.Range("N50").Value = -61248
$ws.Range("H74").Value = 80114
$ws.Range("J74").Value = 80114
$ws.Range("L74").Value = 80114
$ws.Range("N74").Value = -81862
$ws.Range("H77").Value = 80114
$ws.Range("J77").Value = 80114
$ws.Range("L77").Value = 240342
$ws.Range("N77").Value = -249078
$ws.Range("H87").Value = 118760
$ws.Range("J87").Value = 118760
$ws.Range("L87").Value = 118760
$ws.Range("N87").Value = -121132
$ws.Range("H90").Value = 118760
$ws.Range("J90").Value = 118760
$ws.Range("L90").Value = 356280
$ws.Range("N90").Value = -368136
$ws.Range("H105").Value = 52950.715
$ws.Range("I105").Value = 60775.832
$ws.Range("K105").Value = 60775.832
$ws.Range("M105").Value = -59028.832

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1800
$ws.Range("I16").Value = 2133.3333
$ws.Range("K16").Value = 6399.999899999999
$ws.Range("M16").Value = -6226.999899999999
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = $null
$ws.Range("H20").Value = 299
$ws.Range("I20").Value = 299
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 897
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = $null
$ws.Range("M20").Value = -670
$ws.Range("H38").Value = 58.666668
$ws.Range("I38").Value = 33.333332
$ws.Range("K38").Value = 99.999996
$ws.Range("M38").Value = 247.000004
$ws.Range("H99").Value = 6019.125
$ws.Range("I99").Value = 1180.3334
$ws.Range("J99").Value = 8922.4
$ws.Range("K99").Value = 3541.0002
$ws.Range("L99").Value = 26767.2
$ws.Range("M99").Value = -1295.0002
$ws.Range("N99").Value = -31259.2
$ws.Range("H105").Value = 20029
$ws.Range("J105").Value = 20029
$ws.Range("L105").Value = 60087
$ws.Range("N105").Value = -65329
$ws.Range("H121").Value = 707
$ws.Range("J121").Value = 923
$ws.Range("L121").Value = 2769
$ws.Range("N121").Value = -5389
$ws.Range("H129").Value = 3000.2354
$ws.Range("J129").Value = 3283.6667
$ws.Range("L129").Value = 9851.000100000001
$ws.Range("N129").Value = -19851.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4009600
$ws.Range("J21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("N21").Value = -12346
$ws.Range("H30").Value = 4009600
$ws.Range("J30").Value = 12000
$ws.Range("L30").Value = 12000
$ws.Range("N30").Value = -12210
$ws.Range("H57").Value = 11439.777
$ws.Range("J57").Value = 40029
$ws.Range("L57").Value = 40029
$ws.Range("N57").Value = -41669
$ws.Range("H80").Value = 398932.94
$ws.Range("I80").Value = 644890
$ws.Range("J80").Value = 5401.6
$ws.Range("K80").Value = 644890
$ws.Range("L80").Value = 5401.6
$ws.Range("M80").Value = -643892
$ws.Range("N80").Value = -7397.6
$ws.Range("H83").Value = 398932.94
$ws.Range("I83").Value = 644890
$ws.Range("J83").Value = 5401.6
$ws.Range("K83").Value = 3224450
$ws.Range("L83").Value = 27008
$ws.Range("M83").Value = -3219458
$ws.Range("N83").Value = -36992
$ws.Range("H122").Value = 5413.0347
$ws.Range("J122").Value = 6889.875
$ws.Range("L122").Value = 20669.625
$ws.Range("N122").Value = -25569.625
$ws.Range("H132").Value = 58829736
$ws.Range("I132").Value = 76927630
$ws.Range("J132").Value = 11573
$ws.Range("K132").Value = 230782890
$ws.Range("L132").Value = 34719
$ws.Range("M132").Value = -230780360
$ws.Range("N132").Value = -39779

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10859
$ws.Range("I16").Value = 17014.834
$ws.Range("J16").Value = 1625.25
$ws.Range("K16").Value = 17014.834
$ws.Range("L16").Value = 1625.25
$ws.Range("M16").Value = -16844.834
$ws.Range("N16").Value = -1965.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4758
$ws.Range("I81").Value = 3355.3333
$ws.Range("J81").Value = 5599.6
$ws.Range("K81").Value = 6710.6666
$ws.Range("L81").Value = 11199.2
$ws.Range("M81").Value = -5649.6666
$ws.Range("N81").Value = -13321.2
$ws.Range("H84").Value = 4758
$ws.Range("I84").Value = 3355.3333
$ws.Range("J84").Value = 5599.6
$ws.Range("K84").Value = 33553.333
$ws.Range("L84").Value = 55996
$ws.Range("M84").Value = -28249.333
$ws.Range("N84").Value = -66604
